$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update best_known (B) and gap-source (C) values for rows 2-31 (new VNS run results)
$ws.Range("B2").Value = 24038
$ws.Range("C2").Value = 23326
$ws.Range("B3").Value = 24620
$ws.Range("C3").Value = 24114
$ws.Range("B4").Value = 24347
$ws.Range("C4").Value = 23959
$ws.Range("B5").Value = 22913
$ws.Range("C5").Value = 22709
$ws.Range("B6").Value = 24664
$ws.Range("C6").Value = 24219
$ws.Range("B7").Value = 24906
$ws.Range("C7").Value = 25258
$ws.Range("B8").Value = 25313
$ws.Range("C8").Value = 25071
$ws.Range("B9").Value = 23366
$ws.Range("C9").Value = 23337
$ws.Range("B10").Value = 24756
$ws.Range("C10").Value = 24365
$ws.Range("B11").Value = 24096
$ws.Range("C11").Value = 24220
$ws.Range("B12").Value = 42505
$ws.Range("C12").Value = 41351
$ws.Range("B13").Value = 41934
$ws.Range("C13").Value = 40983
$ws.Range("B14").Value = 41737
$ws.Range("C14").Value = 41417
$ws.Range("B15").Value = 44801
$ws.Range("C15").Value = 43992
$ws.Range("B16").Value = 42001
$ws.Range("C16").Value = 40717
$ws.Range("B17").Value = 42946
$ws.Range("C17").Value = 42666
$ws.Range("B18").Value = 41607
$ws.Range("C18").Value = 40661
$ws.Range("B19").Value = 44441
$ws.Range("C19").Value = 44243
$ws.Range("B20").Value = 43372
$ws.Range("C20").Value = 43664
$ws.Range("B21").Value = 44275
$ws.Range("C21").Value = 43837
$ws.Range("B22").Value = 60105
$ws.Range("C22").Value = 59910
$ws.Range("B23").Value = 62231
$ws.Range("C23").Value = 62514
$ws.Range("B24").Value = 59807
$ws.Range("C24").Value = 59708
$ws.Range("B25").Value = 60379
$ws.Range("C25").Value = 60319
$ws.Range("B26").Value = 60700
$ws.Range("C26").Value = 61043
$ws.Range("B27").Value = 59328
$ws.Range("C27").Value = 58889
$ws.Range("B28").Value = 60852
$ws.Range("C28").Value = 60911
$ws.Range("B29").Value = 60409
$ws.Range("C29").Value = 61017
$ws.Range("B30").Value = 59073
$ws.Range("C30").Value = 59146
$ws.Range("B31").Value = 59972
$ws.Range("C31").Value = 60164

# Add newly-available C values for rows 212-214 (previously missing data)
$ws.Range("C212").Value = 54442
$ws.Range("C213").Value = 55460
$ws.Range("C214").Value = 53504

# Update the view state: scroll position + active selection
$ws.Range("F6").Select()

Write-Output "done"
